$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF column (F) values per repull of data / mean calculation
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = -4
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -3
